$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(31).Delete()
